$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these cells keep their text (string) representation instead of
# being auto-converted to numbers/percentages by Excel when the new value
# looks numeric, by pre-formatting the affected cells as Text ("@").
$cellRefs = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "E20", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "E25", "E26", "D27", "E27", "D40", "E40", "D41", "E41", "D42", "E42", "E43", "E44", "D45", "E45", "E46", "E47", "D48", "E48", "E49", "E50")
foreach ($ref in $cellRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated values (kept as text so they match the source data format).
$ws.Range("D2").Value = "260.83"
$ws.Range("E2").Value = "1.42%"
$ws.Range("D3").Value = "27.19"
$ws.Range("E3").Value = "1.65%"
$ws.Range("D4").Value = "4.673"
$ws.Range("E4").Value = "0.82%"
$ws.Range("D5").Value = "0.06179"
$ws.Range("E5").Value = "4.02%"
$ws.Range("D6").Value = "6.666"
$ws.Range("E6").Value = "0.85%"
$ws.Range("D7").Value = "0.8513"
$ws.Range("E7").Value = "-0.65%"
$ws.Range("D8").Value = "0.9178"
$ws.Range("E8").Value = "0.12%"
$ws.Range("D9").Value = "0.1412"
$ws.Range("E9").Value = "2.14%"
$ws.Range("D10").Value = "0.04852"
$ws.Range("E10").Value = "9.92%"
$ws.Range("D11").Value = "0.07086"
$ws.Range("E11").Value = "1.18%"
$ws.Range("D12").Value = "0.03111"
$ws.Range("E12").Value = "3.30%"
$ws.Range("D13").Value = "0.09042"
$ws.Range("E13").Value = "-0.72%"
$ws.Range("D14").Value = "0.001540"
$ws.Range("E14").Value = "0.25%"
$ws.Range("D15").Value = "0.0006147"
$ws.Range("E15").Value = "1.87%"
$ws.Range("D16").Value = "0.006002"
$ws.Range("E16").Value = "-2.07%"
$ws.Range("D17").Value = "3.449"
$ws.Range("E17").Value = "-0.54%"
$ws.Range("D18").Value = "3.154"
$ws.Range("E18").Value = "0.74%"
$ws.Range("D19").Value = "2.179"
$ws.Range("E19").Value = "1.37%"
$ws.Range("E20").Value = "-0.56%"
$ws.Range("E21").Value = "0.24%"
$ws.Range("D22").Value = "4.085"
$ws.Range("E22").Value = "5.72%"
$ws.Range("D23").Value = "0.04236"
$ws.Range("E23").Value = "1.29%"
$ws.Range("D24").Value = "0.001215"
$ws.Range("E24").Value = "0.06%"
$ws.Range("E25").Value = "-15.33%"
$ws.Range("E26").Value = "0.14%"
$ws.Range("D27").Value = "0.0001575"
$ws.Range("E27").Value = "-8.03%"
$ws.Range("D40").Value = "0.03875"
$ws.Range("E40").Value = "1.52%"
$ws.Range("D41").Value = "0.1112"
$ws.Range("E41").Value = "0.55%"
$ws.Range("D42").Value = "0.004098"
$ws.Range("E42").Value = "10.69%"
$ws.Range("E43").Value = "8.17%"
$ws.Range("E44").Value = "-4.67%"
$ws.Range("D45").Value = "0.00005170"
$ws.Range("E45").Value = "2.10%"
$ws.Range("E46").Value = "0.17%"
$ws.Range("E47").Value = "7.97%"
$ws.Range("D48").Value = "0.1624"
$ws.Range("E48").Value = "-32.59%"
$ws.Range("E49").Value = "0.17%"
$ws.Range("E50").Value = "0.17%"
